# Generate Report for Handback
# Update the handback-status report with refreshed timestamps (and the
# zh-cn/de-de "Priority" value for the first file, which moved from a
# human-translation ("ht") hand-off to a machine-translation ("mt") one)
# as of the latest report generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# "Latest HO Xliff Generate Date" for the first file (rows 2 and 3 both
# reference the same value).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-30 20:18:21"
$wsOverview.Range("G3").Value = "2016-08-30 20:18:21"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-30 20:18:16"
$wsZhCn.Range("H3").Value = "2016-08-30 20:18:16"
$wsZhCn.Range("K2").Value = "2016-08-30 20:18:34"
$wsZhCn.Range("K3").Value = "2016-08-30 20:18:34"

# --- de-de sheet ------------------------------------------------------
# Note: de-de's "Correspond Handoff Datetime" for the first file happens to
# share the very same underlying text value as the Overview sheet's
# "Latest HO Xliff Generate Date" ("2016-08-30 20:17:25" -> "...20:18:21"),
# so it must be updated in lock-step with Overview!G2/G3 above.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-30 20:18:21"
$wsDeDe.Range("H3").Value = "2016-08-30 20:18:21"
$wsDeDe.Range("K2").Value = "2016-08-30 20:18:40"
$wsDeDe.Range("K3").Value = "2016-08-30 20:18:40"
